# Updates cryptos list values (price + 1h volume change, and a couple of
# coin-name/link/price row swaps) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "68.291.88"
$ws.Range("E2").Value = "  +1.02%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.348.48"
$ws.Range("E3").Value = "  +0.63%  "

# Row 4: TetherUSD
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5: BNB
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "583.29"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.27%  "

# Row 6: Solana
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "177.16"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.05%  "

# Row 7: USDC
$ws.Range("E7").Value = "  -0.16%  "

# Row 8: XRP
$ws.Range("E8").Value = "  +0.23%  "

# Row 9: Dogecoin
$ws.Range("E9").Value = "  +3.28%  "

# Row 10: Cardano
$ws.Range("E10").Value = "  +0.98%  "

# Row 11: Avalanche
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "47.95"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +5.50%  "

# Row 12: ShibaInu
$ws.Range("E12").Value = "  +1.36%  "

# Row 13: BitcoinCash
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "685.29"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +4.15%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.888.32"

# Row 15: Polkadot
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "8.41"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "68.330.47"
$ws.Range("E16").Value = "  +0.86%  "

# Row 17: TRON
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.119"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.10%  "

# Row 18: WrappedEther
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.331.90"
$ws.Range("E18").Value = "  -0.14%  "

# Row 19: Chainlink
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "17.43"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "

# Row 20: Uniswap
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.18"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.39%  "

# Row 21: Polygon
$ws.Range("E21").Value = "  +0.61%  "

# Row 22: Toncoin
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.44"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "

# Row 23: InternetComputer(DFINITY)
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "16.90"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.59%  "

# Row 24: Litecoin
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "99.98"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.46%  "

# Row 25: PancakeSwap
$ws.Range("E25").Value = "  +1.58%  "

# Row 26: ImmutableX
$ws.Range("E26").Value = "  +1.12%  "

# Row 27: RenderToken
$ws.Range("E27").Value = "  +2.98%  "

# Row 28: EthereumClassic
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "32.97"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.65%  "

# Row 29: Filecoin
$ws.Range("E29").Value = "  +0.90%  "

# Row 30: NEARProtocol
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.94"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -6.05%  "

# Row 31: Bittensor
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "562.27"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.32%  "

# Row 32: Cosmos
$ws.Range("E32").Value = "  +0.76%  "

# Row 33: Hedera
$ws.Range("E33").Value = "  +0.92%  "

# Row 34: OKB
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "57.90"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.26%  "

# Row 35: Dai
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "

# Row 36: Maker
$ws.Range("D36").Value = "3.709.22"
$ws.Range("E36").Value = "  +0.09%  "

# Row 37: dogwifhat
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.30"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.68%  "

# Row 38: Kaspa
$ws.Range("E38").Value = "  +4.12%  "

# Row 39: InjectiveProtocol
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "34.66"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.53%  "

# Row 40: Stacks
$ws.Range("E40").Value = "  +1.69%  "

# Row 41: Fetch.AI
$ws.Range("E41").Value = "  -0.67%  "

# Row 42: PEPE
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0671"
$ws.Range("E42").Value = "  +1.18%  "

# Row 43: TheGraph
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.335"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.72%  "

# Row 44: ApeXProtocol
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.24"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.89%  "

# Row 45: VeChain
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0411"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.28%  "

# Row 46: ThetaToken
$ws.Range("E46").Value = "  +2.38%  "

# Row 47: Stellar
$ws.Range("E47").Value = "  +0.51%  "

# Row 48: FirstDigitalUSD
$ws.Range("E48").Value = "  -0.26%  "

# Row 49: Mantle
$ws.Range("E49").Value = "  -0.29%  "

# Row 50: Monero
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "131.25"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.18%  "

# Row 51: CoreDAO
$ws.Range("E51").Value = "  -0.11%  "
